# Fix bug in total cost: update the "address" column (D) values on the
# "Child" worksheet for rows 2-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Child")

$ws.Range("D2").Value  = "5,-2"
$ws.Range("D3").Value  = "-7,7"
$ws.Range("D4").Value  = "-1,9"
$ws.Range("D5").Value  = "-6,2"
$ws.Range("D6").Value  = "0,9"
$ws.Range("D7").Value  = "-4,1"
$ws.Range("D8").Value  = "4,9"
$ws.Range("D9").Value  = "-1,-3"
$ws.Range("D10").Value = "-2,-9"
$ws.Range("D11").Value = "-9,-10"
$ws.Range("D12").Value = "-3,-7"
$ws.Range("D13").Value = "-1,-8"
$ws.Range("D14").Value = "2,-10"
$ws.Range("D15").Value = "-10,0"
$ws.Range("D16").Value = "-7,9"
$ws.Range("D17").Value = "-1,-6"
$ws.Range("D18").Value = "5,4"
$ws.Range("D19").Value = "2,9"
$ws.Range("D20").Value = "4,8"
$ws.Range("D21").Value = "-5,3"
